$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 38 ("Mise en place de l'environnement de dévellopement") ---
# This shifts everything from row 38 downward by one row (old row 38 -> 39, ... old row 44 -> 45, etc.)
[void]$ws.Rows.Item(38).Insert()

# Bring over the formatting from the row above (37) for the two columns whose
# style differs from the generic body style (A = left border, F = percent style),
# matching what Excel does when a row is inserted in the middle of a formatted table.
[void]$ws.Range("A37").Copy()
[void]$ws.Range("A38").PasteSpecial(-4122)
[void]$ws.Range("F37").Copy()
[void]$ws.Range("F38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new row's data ---
$ws.Range("A38").Value = "Réalisation"
# Set C38 before B38 so the two brand-new shared strings land in the same
# order as the target workbook ("En cours" then "Mise en place ...").
$ws.Range("C38").Value = "En cours"
$ws.Range("B38").Value = "Mise en place de l'environnement de dévellopement"
$ws.Range("D38").Value = 4
$ws.Range("E38").Value = 2.5
$ws.Range("F38").Formula = "=E38/D38"

# --- Grow the AutoFilter range from C1:F44 to C1:F45 (row 38 was inserted
# inside the filtered range, so the filter needs to cover the extra row) ---
$ws.AutoFilterMode = $false
[void]$ws.Range("C1:F45").AutoFilter(1)

# --- Keep the _FilterDatabase defined name in sync with the new filter range ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Planning & Journal'!`$C`$1:`$F`$45"
    }
}

# --- Restore the active selection recorded in the saved workbook ---
[void]$ws.Range("K35").Select()
